$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.072518403766446
$ws.Cells.Item(2, 4).Value = 1.071734474037006
$ws.Cells.Item(2, 5).Value = 1.076252960779337
$ws.Cells.Item(2, 6).Value = 1.086234518589867
$ws.Cells.Item(2, 9).Value = 1.06157989630607
$ws.Cells.Item(2, 10).Value = 1.077437545509292
$ws.Cells.Item(2, 11).Value = 1.074430858838986
$ws.Cells.Item(2, 12).Value = 1.078937380886507
$ws.Cells.Item(2, 13).Value = 1.088892888960647
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.073887812894444
$ws.Cells.Item(3, 4).Value = 1.072835757304792
$ws.Cells.Item(3, 5).Value = 1.077505242709076
$ws.Cells.Item(3, 6).Value = 1.087588028979666
$ws.Cells.Item(3, 9).Value = 1.062110832322892
$ws.Cells.Item(3, 10).Value = 1.078463321328203
$ws.Cells.Item(3, 11).Value = 1.075348049449877
$ws.Cells.Item(3, 12).Value = 1.080006068472967
$ws.Cells.Item(3, 13).Value = 1.090064456341428
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.074772557995499
$ws.Cells.Item(4, 4).Value = 1.073547145343367
$ws.Cells.Item(4, 5).Value = 1.078314015975194
$ws.Cells.Item(4, 6).Value = 1.088462685324938
$ws.Cells.Item(4, 9).Value = 1.062452409196289
$ws.Cells.Item(4, 10).Value = 1.079125210988189
$ws.Cells.Item(4, 11).Value = 1.07593972395369
$ws.Cells.Item(4, 12).Value = 1.080695480968981
$ws.Cells.Item(4, 13).Value = 1.090820836773028
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.075144186420541
$ws.Cells.Item(5, 4).Value = 1.073845925778154
$ws.Cells.Item(5, 5).Value = 1.078653661465731
$ws.Cells.Item(5, 6).Value = 1.088830119571593
$ws.Cells.Item(5, 9).Value = 1.062595538018972
$ws.Cells.Item(5, 10).Value = 1.079403029489043
$ws.Cells.Item(5, 11).Value = 1.076188034622841
$ws.Cells.Item(5, 12).Value = 1.08098481206655
$ws.Cells.Item(5, 13).Value = 1.09113841655431
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.075206565963697
$ws.Cells.Item(6, 4).Value = 1.073896075607804
$ws.Cells.Item(6, 5).Value = 1.078710668312866
$ws.Cells.Item(6, 6).Value = 1.0888917976445
$ws.Cells.Item(6, 9).Value = 1.062619542485075
$ws.Cells.Item(6, 10).Value = 1.079449650768316
$ws.Cells.Item(6, 11).Value = 1.076229702013603
$ws.Cells.Item(6, 12).Value = 1.081033362971675
$ws.Cells.Item(6, 13).Value = 1.091191716133436
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.074777524956682
$ws.Cells.Item(7, 4).Value = 1.073551138785309
$ws.Cells.Item(7, 5).Value = 1.07831855575399
$ws.Cells.Item(7, 6).Value = 1.088467596058661
$ws.Cells.Item(7, 9).Value = 1.062454323534629
$ws.Cells.Item(7, 10).Value = 1.079128924935703
$ws.Cells.Item(7, 11).Value = 1.075943043575092
$ws.Cells.Item(7, 12).Value = 1.080699348975789
$ws.Cells.Item(7, 13).Value = 1.090825081865265
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.072981485868489
$ws.Cells.Item(8, 4).Value = 1.072106912035172
$ws.Cells.Item(8, 5).Value = 1.076676496643899
$ws.Cells.Item(8, 6).Value = 1.086692186202848
$ws.Cells.Item(8, 9).Value = 1.061759738929664
$ws.Cells.Item(8, 10).Value = 1.077784598173006
$ws.Cells.Item(8, 11).Value = 1.074741204604064
$ws.Cells.Item(8, 12).Value = 1.079298986533322
$ws.Cells.Item(8, 13).Value = 1.089289181016534
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.069806004645033
$ws.Cells.Item(9, 4).Value = 1.069552508372835
$ws.Cells.Item(9, 5).Value = 1.07377097695237
$ws.Cells.Item(9, 6).Value = 1.083554599238296
$ws.Cells.Item(9, 9).Value = 1.060520557561419
$ws.Cells.Item(9, 10).Value = 1.075401311360129
$ws.Cells.Item(9, 11).Value = 1.072609382870748
$ws.Cells.Item(9, 12).Value = 1.076815071582798
$ws.Cells.Item(9, 13).Value = 1.086569471428186
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.067681504702727
$ws.Cells.Item(10, 4).Value = 1.06784293774621
$ws.Cells.Item(10, 5).Value = 1.071825569370158
$ws.Cells.Item(10, 6).Value = 1.081456427884367
$ws.Cells.Item(10, 9).Value = 1.059684032647645
$ws.Cells.Item(10, 10).Value = 1.073802486140461
$ws.Cells.Item(10, 11).Value = 1.071178494396455
$ws.Cells.Item(10, 12).Value = 1.075147868078263
$ws.Cells.Item(10, 13).Value = 1.084747108415519
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.066759709130938
$ws.Cells.Item(11, 4).Value = 1.067101044037986
$ws.Cells.Item(11, 5).Value = 1.070981122375452
$ws.Cells.Item(11, 6).Value = 1.080546291927055
$ws.Cells.Item(11, 9).Value = 1.059319304154323
$ws.Cells.Item(11, 10).Value = 1.073107750589969
$ws.Cells.Item(11, 11).Value = 1.070556554809427
$ws.Cells.Item(11, 12).Value = 1.074423211847766
$ws.Cells.Item(11, 13).Value = 1.083955747728346
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.066417024672263
$ws.Cells.Item(12, 4).Value = 1.066825220532879
$ws.Cells.Item(12, 5).Value = 1.070667139389561
$ws.Cells.Item(12, 6).Value = 1.080207977841417
$ws.Cells.Item(12, 9).Value = 1.05918344780693
$ws.Cells.Item(12, 10).Value = 1.07284932377171
$ws.Cells.Item(12, 11).Value = 1.070325180092357
$ws.Cells.Item(12, 12).Value = 1.074153623743958
$ws.Cells.Item(12, 13).Value = 1.083661454673407
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.066490544823996
$ws.Cells.Item(13, 4).Value = 1.066884397056838
$ws.Cells.Item(13, 5).Value = 1.070734504288631
$ws.Cells.Item(13, 6).Value = 1.08028055873841
$ws.Cells.Item(13, 9).Value = 1.059212606699472
$ws.Cells.Item(13, 10).Value = 1.072904774082825
$ws.Cells.Item(13, 11).Value = 1.070374827062022
$ws.Cells.Item(13, 12).Value = 1.074211470341856
$ws.Cells.Item(13, 13).Value = 1.08372459730911
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.066731388654238
$ws.Cells.Item(14, 4).Value = 1.067078249533466
$ws.Cells.Item(14, 5).Value = 1.070955174963404
$ws.Cells.Item(14, 6).Value = 1.080518331898399
$ws.Cells.Item(14, 9).Value = 1.059308081995592
$ws.Cells.Item(14, 10).Value = 1.073086396550685
$ws.Cells.Item(14, 11).Value = 1.070537436647543
$ws.Cells.Item(14, 12).Value = 1.074400936183084
$ws.Cells.Item(14, 13).Value = 1.083931428475725
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.066879742038351
$ws.Cells.Item(15, 4).Value = 1.067197655012168
$ws.Cells.Item(15, 5).Value = 1.071091095149668
$ws.Cells.Item(15, 6).Value = 1.080664798587898
$ws.Cells.Item(15, 9).Value = 1.059366857034231
$ws.Cells.Item(15, 10).Value = 1.07319825080861
$ws.Cells.Item(15, 11).Value = 1.070637578109974
$ws.Cells.Item(15, 12).Value = 1.074517616705372
$ws.Cells.Item(15, 13).Value = 1.08405881789128
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.067742641167735
$ws.Cells.Item(16, 4).Value = 1.067892139843755
$ws.Cells.Item(16, 5).Value = 1.071881568323864
$ws.Cells.Item(16, 6).Value = 1.081516796045176
$ws.Cells.Item(16, 9).Value = 1.059708185370366
$ws.Cells.Item(16, 10).Value = 1.073848541710122
$ws.Cells.Item(16, 11).Value = 1.071219720444573
$ws.Cells.Item(16, 12).Value = 1.075195902761954
$ws.Cells.Item(16, 13).Value = 1.084799580167481
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.068283408497727
$ws.Cells.Item(17, 4).Value = 1.068327329676091
$ws.Cells.Item(17, 5).Value = 1.072376851859118
$ws.Cells.Item(17, 6).Value = 1.082050794807497
$ws.Cells.Item(17, 9).Value = 1.05992161799101
$ws.Cells.Item(17, 10).Value = 1.074255796597578
$ws.Cells.Item(17, 11).Value = 1.071584248966395
$ws.Cells.Item(17, 12).Value = 1.075620634289883
$ws.Cells.Item(17, 13).Value = 1.085263629947494
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.068598648692313
$ws.Cells.Item(18, 4).Value = 1.068581010806738
$ws.Cells.Item(18, 5).Value = 1.072665542816394
$ws.Cells.Item(18, 6).Value = 1.082362111980109
$ws.Cells.Item(18, 9).Value = 1.060045867889656
$ws.Cells.Item(18, 10).Value = 1.074493107040401
$ws.Cells.Item(18, 11).Value = 1.071796645346099
$ws.Cells.Item(18, 12).Value = 1.075868108448335
$ws.Cells.Item(18, 13).Value = 1.085534084162106
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.068706107095043
$ws.Cells.Item(19, 4).Value = 1.068667482912433
$ws.Cells.Item(19, 5).Value = 1.072763945317337
$ws.Cells.Item(19, 6).Value = 1.082468236952684
$ws.Cells.Item(19, 9).Value = 1.060088192987352
$ws.Cells.Item(19, 10).Value = 1.074573984167787
$ws.Cells.Item(19, 11).Value = 1.071869028673827
$ws.Cells.Item(19, 12).Value = 1.075952446046704
$ws.Cells.Item(19, 13).Value = 1.085626265238113
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.068225407977327
$ws.Cells.Item(20, 4).Value = 1.068280654278549
$ws.Cells.Item(20, 5).Value = 1.072323733311102
$ws.Cells.Item(20, 6).Value = 1.081993517903728
$ws.Cells.Item(20, 9).Value = 1.059898743718733
$ws.Cells.Item(20, 10).Value = 1.074212126292803
$ws.Cells.Item(20, 11).Value = 1.071545162015133
$ws.Cells.Item(20, 12).Value = 1.075575092025228
$ws.Cells.Item(20, 13).Value = 1.085213864409867
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.066660474171745
$ws.Cells.Item(21, 4).Value = 1.067021171784081
$ws.Cells.Item(21, 5).Value = 1.070890201773277
$ws.Cells.Item(21, 6).Value = 1.080448320557786
$ws.Cells.Item(21, 9).Value = 1.05927997741466
$ws.Cells.Item(21, 10).Value = 1.073032923549054
$ws.Cells.Item(21, 11).Value = 1.070489562116988
$ws.Cells.Item(21, 12).Value = 1.074345154816121
$ws.Cells.Item(21, 13).Value = 1.08387053146547
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.065674865204817
$ws.Cells.Item(22, 4).Value = 1.06622782957395
$ws.Cells.Item(22, 5).Value = 1.06998704187309
$ws.Cells.Item(22, 6).Value = 1.079475349333226
$ws.Cells.Item(22, 9).Value = 1.058888734635826
$ws.Cells.Item(22, 10).Value = 1.07228936216628
$ws.Cells.Item(22, 11).Value = 1.069823787059049
$ws.Cells.Item(22, 12).Value = 1.073569420110345
$ws.Cells.Item(22, 13).Value = 1.08302391716429
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.066197516137464
$ws.Cells.Item(23, 4).Value = 1.066648534942496
$ws.Cells.Item(23, 5).Value = 1.070466000738366
$ws.Cells.Item(23, 6).Value = 1.079991278906503
$ws.Cells.Item(23, 9).Value = 1.059096349427042
$ws.Cells.Item(23, 10).Value = 1.072683743736633
$ws.Cells.Item(23, 11).Value = 1.070176925543544
$ws.Cells.Item(23, 12).Value = 1.073980883550571
$ws.Cells.Item(23, 13).Value = 1.083472915763462
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.068251616472294
$ws.Cells.Item(24, 4).Value = 1.068301745368639
$ws.Cells.Item(24, 5).Value = 1.072347735912627
$ws.Cells.Item(24, 6).Value = 1.082019399350401
$ws.Cells.Item(24, 9).Value = 1.059909080364468
$ws.Cells.Item(24, 10).Value = 1.07423185974777
$ws.Cells.Item(24, 11).Value = 1.071562824426808
$ws.Cells.Item(24, 12).Value = 1.075595671431293
$ws.Cells.Item(24, 13).Value = 1.085236351987877
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.070628239895182
$ws.Cells.Item(25, 4).Value = 1.070214033709292
$ws.Cells.Item(25, 5).Value = 1.074523579055632
$ws.Cells.Item(25, 6).Value = 1.084366852156094
$ws.Cells.Item(25, 9).Value = 1.060842737544398
$ws.Cells.Item(25, 10).Value = 1.07601918444925
$ws.Cells.Item(25, 11).Value = 1.073162197121180
$ws.Cells.Item(25, 12).Value = 1.077459185516622
$ws.Cells.Item(25, 13).Value = 1.087274184320544
